$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-obsolete low_bound/high_bound/step_plus cells on rows 2 and 3
# (full Clear so the cells drop out of the sheet entirely, not just their value)
$ws.Range("H2:J2").Clear()
$ws.Range("H3:I3").Clear()

# Update D7 from 1 to 0
$ws.Range("D7").Value = 0

# Update the saved selection to E26
$ws.Range("E26").Select()
